$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New columns H (credit_purpose) and I (credit_purpose_ff) -------------

# Column widths for the two new columns (raw stored widths 15.5 / 16.5;
# the COM ColumnWidth setter applies a +5/6 offset internally, so we
# back that out here to land on the exact target stored width).
$ws.Columns.Item(8).ColumnWidth = 14.666666666666666
$ws.Columns.Item(9).ColumnWidth = 15.666666666666666

# Header row height shrinks back to the normal 17 now that it's a single
# line header again.
$ws.Rows.Item(1).RowHeight = 17

# Headers
$ws.Range("H1").WrapText = $true
$ws.Range("H1").Value = "credit_purpose"
$ws.Range("I1").WrapText = $true
$ws.Range("I1").Value = "credit_purpose_ff"

# Row 2
$ws.Range("H2").WrapText = $true
$ws.Range("H2").Value = "1;2;3"
$ws.Range("I2").WrapText = $true
$ws.Range("I2").Value = "abc;def"

# Row 3
$ws.Range("H3").WrapText = $true
$ws.Range("H3").Value = 999
$ws.Range("I3").WrapText = $true
$ws.Range("I3").Value = "abc;def"

# Row 4
$ws.Range("H4").WrapText = $true
$ws.Range("H4").Value = "1;2;3;4;5;6;7;8"
$ws.Range("I4").WrapText = $true
$ws.Range("I4").Value = "abc;def"

# Row 5
$ws.Range("H5").WrapText = $true
$ws.Range("H5").Value = "10000;1200;1;2"
$ws.Range("I5").WrapText = $true
$ws.Range("I5").Value = "abc;def"

# Row 6
$ws.Range("H6").WrapText = $true
$ws.Range("H6").Value = 988
$ws.Range("I6").WrapText = $true
$ws.Range("I6").Value = "abc;def"

# Row 7
$ws.Range("H7").WrapText = $true
$ws.Range("H7").Value = "1;1"

# Row 8
$ws.Range("H8").WrapText = $true
$ws.Range("H8").Value = 977

# Row 9 (H9 first; I9's new string "def" is introduced after H11's "1;988"
# below, to match the order new entries were appended to the shared string
# table)
$ws.Range("H9").WrapText = $true
$ws.Range("H9").Value = 977

# Row 10 - also update the multi_value_field_restriction value in D10
$ws.Range("D10").Value = 9
$ws.Range("H10").WrapText = $true
$ws.Range("H10").Value = "977;1"
$ws.Range("I10").WrapText = $true
$ws.Range("I10").Value = "abc;def"

# Row 11 - also update the multi_value_field_restriction value in D11
$ws.Range("D11").Value = 10
$ws.Range("H11").WrapText = $true
$ws.Range("H11").Value = "1;988"
$ws.Range("I11").WrapText = $true
$ws.Range("I11").Value = "abc;def"

# Back to row 9's I cell now that "1;988" has been registered
$ws.Range("I9").WrapText = $true
$ws.Range("I9").Value = "def"

# --- View state -------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("I10").Select()
